# 1.1.9 re-corrected once again the amounts and taxes being entered into BOSS
#
# - Marks the application/workbook window as minimized.
# - Inserts a new checklist row on Sheet2 (row 15) recording the new
#   "Branded for Wellbridge" task (dated 10/15/2025), pushing the existing
#   rows down by one.
# - Restores the cell selection that was active on Sheet2 when the file
#   was saved.

$wb = $excel.ActiveWorkbook

# The workbook window is minimized when this checklist is saved.
try {
    $excel.WindowState = -4140   # xlMinimized
} catch {
}
try {
    $wb.Windows.Item(1).WindowState = -4140   # xlMinimized
} catch {
}
try {
    $excel.ActiveWindow.WindowState = -4140   # xlMinimized
} catch {
}

$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Insert a brand new row above the existing row 15 ("Adjust for Cell
# Phones" / 10/21/2025), shifting it (and everything below it) down by one.
$ws.Rows("15").Insert()

# Fill in the new checklist entry.
$ws.Range("A15").Value = 45945
$ws.Range("B15").Value = "Done"
$ws.Range("C15").Value = "Branded for Wellbridge"

# Restore the saved selection.
$ws.Range("I16").Select()
